# Edit: add a new "2022-Q4" worksheet (as 2nd sheet) with fund-holdings data,
# and insert a corresponding summary row at the top of the "总计" (total) sheet's
# data table, renumbering the existing running index.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the "2022-Q4" worksheet by cloning "2022-Q3" (so it inherits the
#    same column layout / header style / text-vs-number cell typing / page
#    margins) and inserting it right after "总计" (i.e. before the current
#    2nd sheet).
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q3")
$template.Copy($wb.Worksheets.Item(2))
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

$q4Data = @(
    @("012428", "华夏核心制造混合A", "32.52", "76.47", "4.00", "1.3008", 8),
    @("013389", "华夏成长先锋一年持有混合A", "11.63", "80.24", "5.47", "0.6362", 6),
    @("012429", "华夏核心制造混合C", "9.69", "76.47", "4.00", "0.3876", 8),
    @("013390", "华夏成长先锋一年持有混合C", "3.59", "80.24", "5.47", "0.1964", 6),
    @("005434", "鹏华睿投灵活配置混合A", "6.02", "82.78", "3.10", "0.1866", 2),
    @("014410", "华夏时代领航两年持有混合A", "2.51", "70.90", "4.14", "0.1039", 7),
    @("014411", "华夏时代领航两年持有混合C", "0.45", "70.90", "4.14", "0.0186", 7),
    @("016950", "鹏华睿投灵活配置混合C", "0.16", "82.78", "3.10", "0.0050", 2)
)

# The template sheet only had 4 data rows (rows 2-5); extend the "A" column's
# running-index style down to rows 6-9 to match it before filling values in.
$newSheet.Cells.Item(2, 1).Copy()
$newSheet.Range("A6:A9").PasteSpecial(-4122)
$newSheet.Application.CutCopyMode = $false

function Set-TextValue($cell, $val) {
    # Force the cell to accept a numeric-looking string as literal text
    # (matching the source data's inlineStr cell type), then drop back to
    # the default/unstyled cell format so no stray explicit style lingers.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

for ($i = 0; $i -lt $q4Data.Count; $i++) {
    $row = $i + 2
    $rec = $q4Data[$i]

    $newSheet.Cells.Item($row, 1).Value = $i

    Set-TextValue $newSheet.Cells.Item($row, 2) $rec[0]
    Set-TextValue $newSheet.Cells.Item($row, 3) $rec[1]
    Set-TextValue $newSheet.Cells.Item($row, 4) $rec[2]
    Set-TextValue $newSheet.Cells.Item($row, 5) $rec[3]
    Set-TextValue $newSheet.Cells.Item($row, 6) $rec[4]
    Set-TextValue $newSheet.Cells.Item($row, 7) $rec[5]

    $newSheet.Cells.Item($row, 8).Value = $rec[6]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row right under the header
#    for "2022-Q4", pushing the existing quarters down by one row, and
#    renumber the leading running-index column (A) to stay sequential.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert(-4121)
$total.Range("B2:D2").ClearFormats()

$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)
$total.Application.CutCopyMode = $false

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 8
$total.Cells.Item(2, 4).Value = 2.84

for ($r = 3; $r -le 10; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------------
# 3. Restore the originally-active sheet/tab ("2020-Q4" was the selected tab
#    in the source workbook).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
